$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 416, shifting existing rows 416:498 down to 417:499
$ws.Rows("416:416").Insert()

# Populate the newly inserted row 416 with the new data record
$ws.Range("A416").Value = 9
$ws.Range("B416").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C416").Value = "Metropolitana"
$ws.Range("D416").Value = 45209
$ws.Range("E416").Value = 13
$ws.Range("F416").Value = 300000001
$ws.Range("G416").Value = "Rabanito"
$ws.Range("H416").Value = "Sin especificar"
$ws.Range("I416").Value = "Primera"
$ws.Range("J416").Value = 7000
$ws.Range("K416").Value = 3000
$ws.Range("L416").Value = 3000
$ws.Range("M416").Value = 3000
$ws.Range("N416").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O416").Value = "Provincia de Chacabuco"
$ws.Range("P416").Value = 30
$ws.Range("Q416").Value = 100
$ws.Range("R416").Value = "Hortaliza"
